# MassWateR Sites Template - add an "Instructions" worksheet after "Sites"
# describing each field of the Sites tab, matching the upstream commit
# "update all templates and sample files in inst #48".

$wb = $excel.ActiveWorkbook
$sites = $wb.Worksheets.Item("Sites")

# ---------------------------------------------------------------------
# 1. Create the new sheet right after "Sites"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Add([Type]::Missing, $sites)
$ws.Name = "Instructions"

# ---------------------------------------------------------------------
# 2. Column widths (approximate character widths used by the template)
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 29.5
$ws.Columns.Item(2).ColumnWidth = 95.66666667
$ws.Columns.Item(3).ColumnWidth = 20.66666667
$ws.Columns.Item(4).ColumnWidth = 20.66666667
$ws.Columns.Item(5).ColumnWidth = 24.16666667

# ---------------------------------------------------------------------
# 3. Free-form intro text (rows 1-2)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "The Sites tab must be formatted exactly like the Sites template, with all of the following fields."
$ws.Range("A2").Value = "The Sites tab must be the first tab in this workbook."
$ws.Range("A1:A2").Font.Bold = $true
$ws.Range("A1:A2").Font.Color = 12611584

$ws.Range("C1").Value = "Template updated 5/19/23"
$ws.Range("C1").Font.Color = 255

# ---------------------------------------------------------------------
# 4. Header row (row 4) - bold-ish grid header, centered, boxed with a
#    thick bottom rule
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "Field"
$ws.Range("B4").Value = "Instructions"
$ws.Range("C4").Value = "Example"
$ws.Range("D4").Value = "Available Values"
$ws.Range("E4").Value = "Required?"

$ws.Range("A4:E4").Borders.LineStyle = 1
$ws.Range("A4:E4").Borders(9).Weight = -4138
$ws.Range("A4:E4").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 5. Data rows (5-9) - one row per Sites field
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "Monitoring Location ID"
$ws.Range("B5").Value = "Location ID that is used in your Results file.  Must match exactly."
$ws.Range("C5").Value = "ABT-010"
$ws.Range("D5").Value = "any"
$ws.Range("E5").Value = "Required"

$ws.Range("A6").Value = "Monitoring Location Name"
$ws.Range("B6").Value = "Name of monitoring location."
$ws.Range("C6").Value = "477 Lowell Rd, Concord"
$ws.Range("D6").Value = "any"
$ws.Range("E6").Value = "Required for WQX"

$ws.Range("A7").Value = "Monitoring Location Latitude "
$ws.Range("B7").Value = "Latitude of monitoring location in decimal form.  At least 5 decimals."
$ws.Range("C7").Value = 42.470370000000003
$ws.Range("D7").Value = "any"
$ws.Range("E7").Value = "Required for mapping"

$ws.Range("A8").Value = "Monitoring Location Longitude"
$ws.Range("B8").Value = "Longitude of monitoring location in decimal form.  At least 5 decimals."
$ws.Range("C8").Value = -71.362578999999997
$ws.Range("D8").Value = "any"
$ws.Range("E8").Value = "Required for mapping"

$ws.Rows.Item(9).RowHeight = 30
$ws.Range("A9").Value = "Location Group"
$ws.Range("B9").Value = "An optional free-form grouping attribute.  This will allow you to summarize locations by group in the graphing and mapping analysis functions."
$ws.Range("C9").Value = "Lower Assabet"
$ws.Range("D9").Value = "any"
$ws.Range("E9").Value = "Optional"

# Empty, bold-styled spacer cell below the table
$ws.Range("B11").Font.Bold = $true

# ---------------------------------------------------------------------
# 6. Formatting for the data block (rows 5-9)
#    * Thin box border all around every cell
#    * No top border on row 5 (sits right below the header's rule)
#    * Top-aligned everywhere, wrapped text in column B, centered in C:E
#    * Italic font for the "Available Values" / "Required?" columns
# ---------------------------------------------------------------------
$ws.Range("A5:E9").Borders.LineStyle = 1
$ws.Range("A5:E5").Borders(8).LineStyle = 0
$ws.Range("A5:E9").VerticalAlignment = -4160
$ws.Range("C5:E9").HorizontalAlignment = -4108
$ws.Range("B5:B9").WrapText = $true
$ws.Range("D5:E9").Font.Italic = $true

# ---------------------------------------------------------------------
# 7. Freeze panes (column A + rows 1-4 frozen) and initial selection
# ---------------------------------------------------------------------
$ws.Range("B5").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B5").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C2").Select()

# ---------------------------------------------------------------------
# 8. Restore "Sites" as the active sheet/selection
# ---------------------------------------------------------------------
$sites.Activate()
$sites.Range("E2").Select()
